$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.639.55"
$ws.Range("E2").Value = "  +5.36%  "
$ws.Range("D3").Value = "2.233.02"
$ws.Range("E3").Value = "  +3.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.57"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.620"
$ws.Range("E6").Value = "  -2.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.70"
$ws.Range("E7").Value = "  -2.73%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.402"
$ws.Range("E9").Value = "  +1.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.86"
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0876"
$ws.Range("E11").Value = "  +2.96%  "
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "2.559.51"
$ws.Range("E13").Value = "  +2.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.65"
$ws.Range("E14").Value = "  -2.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.02"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.798"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.57"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").Value = "2.234.62"
$ws.Range("E18").Value = "  +3.40%  "
$ws.Range("D19").Value = "41.497.80"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.23"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("D21").Value = "0.0₃0903"
$ws.Range("E21").Value = "  +6.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.00"
$ws.Range("E22").Value = "  -3.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "247.55"
$ws.Range("E23").Value = "  +7.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.39"
$ws.Range("E25").Value = "  +3.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.36"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.57"
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.99"
$ws.Range("E28").Value = "  -1.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.142"
$ws.Range("E29").Value = "  +2.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.09"
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.43"
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.79"
$ws.Range("E32").Value = "  +5.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.122"
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.95"
$ws.Range("E34").Value = "  +5.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.61"
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0624"
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.77"
$ws.Range("E37").Value = "  +2.71%  "
$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.65"
$ws.Range("E38").Value = "  -5.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.38"
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.997"
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("B41").Value = "TerraClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.000235"
$ws.Range("E41").Value = "  +28.27%  "
$ws.Range("B42").Value = "FTXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.83"
$ws.Range("E42").Value = "  +4.26%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0236"
$ws.Range("E43").Value = "  +3.94%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.53"
$ws.Range("E44").Value = "  +8.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.74"
$ws.Range("E45").Value = "  -2.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0961"
$ws.Range("E46").Value = "  +4.03%  "
$ws.Range("D47").Value = "1.488.42"
$ws.Range("E47").Value = "  -2.30%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.78"
$ws.Range("E48").Value = "  -5.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.18"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.78"
$ws.Range("E50").Value = "  -1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.08"
$ws.Range("E51").Value = "  -2.30%  "
